$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header updates
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 updates
$ws.Range("B2").Value = 4.0460026821195827
$ws.Range("C2").Value = $null
$ws.Range("D2").Value = 2.6298241882695685
$ws.Range("E2").Value = $null

# Row 3 updates
$ws.Range("B3").Value = 3.8899910489589478
$ws.Range("C3").Value = 5.3604036204009624
$ws.Range("D3").Value = 2.7694101486199005
$ws.Range("E3").Value = 6.8392903583500395

# Update selection range to reflect the new selection B1:E3
$ws.Range("B1:E3").Select()
